$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3
$ws.Range("U3").Value = 110
$ws.Range("V3").Value = 1024
$ws.Range("W3").Value = 1147
$ws.Range("X3").Value = 110
$ws.Range("Y3").Value = 210
$ws.Range("Z3").Value = 260
$ws.Range("AA3").Value = 571
$ws.Range("AB3").Value = 640
$ws.Range("AC3").Value = 1024
$ws.Range("AD3").Value = 1147
$ws.Range("AE3").Value = 1440
$ws.Range("AF3").Value = 1480
$ws.Range("AG3").Value = 361
$ws.Range("AH3").Value = 861
$ws.Range("AI3").Value = 0.4192799070847851
$ws.Range("AK3").Value = 36.6

# Row 4
$ws.Range("Y4").Value = 160
$ws.Range("Z4").Value = 240
$ws.Range("AA4").Value = 539
$ws.Range("AB4").Value = 720
$ws.Range("AE4").Value = 1290
$ws.Range("AF4").Value = 1370
$ws.Range("AG4").Value = 379
$ws.Range("AH4").Value = 754
$ws.Range("AI4").Value = 0.5026525198938993

# Row 10
$ws.Range("Y10").Value = 840
$ws.Range("Z10").Value = 260
$ws.Range("AA10").Value = 455
$ws.Range("AB10").Value = 550
$ws.Range("AE10").Value = 1660
$ws.Range("AF10").Value = 1080
$ws.Range("AG10").Value = -384.9999999999999
$ws.Range("AH10").Value = 1202
$ws.Range("AI10").Value = -0.3202995008319467
$ws.Range("AK10").Value = 137.5

# Row 15
$ws.Range("AE15").Value = 1503
$ws.Range("AF15").Value = 1586
$ws.Range("AK15").Value = 39.7

# Row 23
$ws.Range("Y23").Value = 1003
$ws.Range("Z23").Value = 281
$ws.Range("AA23").Value = 561
$ws.Range("AB23").Value = 606
$ws.Range("AE23").Value = 2244
$ws.Range("AF23").Value = 1522
$ws.Range("AG23").Value = -442.0000000000001
$ws.Range("AK23").Value = 57.5
